$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.488.11"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.50%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'3.757.70"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -1.83%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("E4").Value = "'  -0.11%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'594.50"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.82%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'169.18"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +1.66%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'3.757.42"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -1.73%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("E8").Value = "'  -0.14%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("E9").Value = "'  -0.18%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'0.164"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +0.48%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'6.48"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +0.78%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("E12").Value = "'  -0.56%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("E13").Value = "'  +4.71%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'36.63"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -0.32%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'4.387.79"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -1.87%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'3.757.76"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -1.93%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'18.80"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +4.03%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'67.448.46"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.86%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("E19").Value = "'  -2.01%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("E20").Value = "'  +0.89%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'10.52"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -3.26%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'468.27"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +0.91%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'0.721"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -0.96%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'0.0000148"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -6.85%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'83.81"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +1.25%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("E26").Value = "'  +0.55%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = "'12.17"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +1.22%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("E28").Value = "'  +2.99%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("E29").Value = "'  +0.26%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("E30").Value = "'  -1.93%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = "'3.905.43"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -1.84%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("D32").Value = "'7.67"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +1.67%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("E33").Value = "'  -2.01%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").Value = "'30.41"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -1.79%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").Value = "'9.13"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -4.09%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("D36").Value = "'3.721.51"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -1.79%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").Value = "'3.82"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +7.13%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("E38").Value = "'  +0.81%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = "'5.89"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +0.26%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "'0.138"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -1.76%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = "'0.998"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -2.16%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = "'0.998"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -0.26%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("E43").Value = "'  +0.84%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("D45").Value = "'8.70"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +1.26%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("E46").Value = "'  -0.89%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = "'45.80"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -2.44%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'396.72"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -5.31%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("B49").Value = "'FLOKI"
$ws.Range("B49").Style = "Normal"
$ws.Range("C49").Value = "'https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("C49").Style = "Normal"
$ws.Range("D49").Value = "'0.000269"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -8.90%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("B50").Value = "'Monero"
$ws.Range("B50").Style = "Normal"
$ws.Range("C50").Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("C50").Style = "Normal"
$ws.Range("D50").Value = "'141.41"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -0.14%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Value = "'0.0353"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -1.07%  "
$ws.Range("E51").Style = "Normal"
